# Fix elective course scheduling to use common time slots for both sections.
# Updates the Section_A and Section_B timetables so CS151 (Elective) lands on
# the same day/slot (Thu, 15:30-17:00) for both sections, and reshuffles the
# other course placements accordingly.

$wb = $excel.ActiveWorkbook

# ---- Section_A ----
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("D2").Value = "MA101"
$wsA.Range("E2").Value = "Free"

$wsA.Range("B3").Value = "DS101"
$wsA.Range("C3").Value = "Free"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "HS101"
$wsA.Range("C5").Value = "HS101"
$wsA.Range("D5").Value = "EC101"
$wsA.Range("E5").Value = "Free"

$wsA.Range("B6").Value = "MA102"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS151 (Elective)"
$wsA.Range("F6").Value = "DS101"

$wsA.Range("B7").Value = "EC101"
$wsA.Range("C7").Value = "CS101"
$wsA.Range("D7").Value = "HS101"
$wsA.Range("E7").Value = "MA102"

# ---- Section_B ----
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "MA101"
$wsB.Range("C2").Value = "HS101"
$wsB.Range("D2").Value = "DS101"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "CS101"

$wsB.Range("B3").Value = "EC101"
$wsB.Range("E3").Value = "MA101"
$wsB.Range("F3").Value = "HS101"

$wsB.Range("B5").Value = "MA102"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "HS101"
$wsB.Range("F5").Value = "Free"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "DS101"
$wsB.Range("D6").Value = "CS101"
$wsB.Range("E6").Value = "CS151 (Elective)"

$wsB.Range("C7").Value = "EC101"
$wsB.Range("D7").Value = "EC101"
$wsB.Range("E7").Value = "CS101"
